$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded for the zh-cn and de-de report rows.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 00:22:06"
$wsZhCn.Range("E5").Value = "2016-03-24 00:22:06"
$wsZhCn.Range("H3").Value = "2016-03-24 00:22:31"
$wsZhCn.Range("H5").Value = "2016-03-24 00:22:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 00:22:11"
$wsDeDe.Range("E5").Value = "2016-03-24 00:22:11"
$wsDeDe.Range("H3").Value = "2016-03-24 00:22:38"
$wsDeDe.Range("H5").Value = "2016-03-24 00:22:38"
